$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) now reflect a completed
# handback that is in sync with en-US instead of merely "ready for handoff".
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn/de-de status columns so the longer status text fits.
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: report the handback - target file / handback file / handback
# datetime columns for both rows are now populated.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$mdAddr1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8813437210004909a83c2c7ffc42f08730bb37a/e2e/8d499f46-2b67-4fd8-b11e-648ab1713868.md"
$mdAddr2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8813437210004909a83c2c7ffc42f08730bb37a/e2e/d8893950-1f4a-4bea-9f77-eedafdb26f91.md"

$zhcn.Range("I2").Value = "8d499f46-2b67-4fd8-b11e-648ab1713868.md"
$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdAddr1, "", "", "8d499f46-2b67-4fd8-b11e-648ab1713868.md")
$zhcn.Range("J2").Value = "8d499f46-2b67-4fd8-b11e-648ab1713868.aad6247552f073641def449800e13234aab0d1f0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-28 20:30:17"

$zhcn.Range("I3").Value = "d8893950-1f4a-4bea-9f77-eedafdb26f91.md"
$zhcn.Range("I3").Style = "HyperLink"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdAddr2, "", "", "d8893950-1f4a-4bea-9f77-eedafdb26f91.md")
$zhcn.Range("J3").Value = "d8893950-1f4a-4bea-9f77-eedafdb26f91.c7231397a7290c9f66aaad1616e06ce4643f83a4.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-28 20:30:17"

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, different handback timestamp.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("I2").Value = "8d499f46-2b67-4fd8-b11e-648ab1713868.md"
$dede.Range("I2").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("I2"), $mdAddr1, "", "", "8d499f46-2b67-4fd8-b11e-648ab1713868.md")
$dede.Range("J2").Value = "8d499f46-2b67-4fd8-b11e-648ab1713868.aad6247552f073641def449800e13234aab0d1f0.de-de.xlf"
$dede.Range("K2").Value = "2016-08-28 20:30:23"

$dede.Range("I3").Value = "d8893950-1f4a-4bea-9f77-eedafdb26f91.md"
$dede.Range("I3").Style = "HyperLink"
$dede.Hyperlinks.Add($dede.Range("I3"), $mdAddr2, "", "", "d8893950-1f4a-4bea-9f77-eedafdb26f91.md")
$dede.Range("J3").Value = "d8893950-1f4a-4bea-9f77-eedafdb26f91.c7231397a7290c9f66aaad1616e06ce4643f83a4.de-de.xlf"
$dede.Range("K3").Value = "2016-08-28 20:30:23"

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

Write-Host "Handback report generated"
